$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.550.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").Value = "'3.611.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'202.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.34%  "

$ws.Range("D6").Value = "'597.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.68%  "

$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +7.18%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "'53.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").Value = "'0.0000302"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "'9.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "'4.185.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "

$ws.Range("D15").Value = "'677.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +14.38%  "

$ws.Range("D16").Value = "'70.617.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("B17").Value = "'Chainlink"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'19.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("B18").Value = "'Uniswap"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("B19").Value = "'WrappedEther"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.584.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("E21").Value = "  +1.87%  "

$ws.Range("D22").Value = "'18.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.08%  "

$ws.Range("D23").Value = "'110.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.05%  "

$ws.Range("E24").Value = "  +3.69%  "

$ws.Range("D25").Value = "'4.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("D26").Value = "'3.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").Value = "'10.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").Value = "'10.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.53%  "

$ws.Range("D30").Value = "'34.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.57%  "

$ws.Range("D31").Value = "'4.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.47%  "

$ws.Range("D32").Value = "'7.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("D33").Value = "'12.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.40%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'63.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "

$ws.Range("D36").Value = "'0.0₃0854"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.47%  "

$ws.Range("D37").Value = "'3.879.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "'513.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("E40").Value = "  -5.13%  "

$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("D42").Value = "'36.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.92%  "

$ws.Range("D43").Value = "'0.385"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.84%  "

$ws.Range("E44").Value = "  +3.46%  "

$ws.Range("D45").Value = "'0.0468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.84%  "

$ws.Range("E46").Value = "  +10.19%  "

$ws.Range("E47").Value = "  +1.46%  "

$ws.Range("D48").Value = "'0.141"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("D49").Value = "'8.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("E51").Value = "  +23.26%  "
